$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (2026-02-21 -> 2026-02-22, i.e. 46074 -> 46075) for every data row.
$ws.Range("C2:C255").Value = 46075
